# Modifiche login lato server
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4: was a bare "Creazione Login" note; becomes a full task row describing
# the fake-login work, with dates + a longer description/notes pair.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 42726
$ws.Range("C4").Value = 42727
$ws.Range("E4").Value = "Create le chiamate principali per l'autenticazione fake."
$ws.Range("F4").Value = "Bisogna trovare un modulo per fare un'autenticazione sicura. Bisogna agganciare un DB. Bisogna creare anche un sistema di registrazione al sito"
$ws.Range("D4").Value = "Creazione Login FAKE"
$ws.Range("B4:C4").NumberFormat = "mm-dd-yy"
$ws.Rows.Item(4).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 5: "Gestione Routing" gains a Descrizione/Note pair.
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "Se l'utente esiste bisogna redirezionare alla home page"
$ws.Range("F5").Value = "Routing - passaggio valore a componenti figli"
$ws.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------------
# Selection moved to F3 (from D9) in the refreshed sheet.
# ---------------------------------------------------------------------------
$ws.Range("F3").Select()

Write-Output "edit applied"
